# Generate Report for Handoff
# Adds two new tracked files (cbb11c9f-... and f235f8b0-...) to each of the
# three report sheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------------

$overview.Range("A4").Value = "cbb11c9f-9aa9-40a5-b27b-62f0134d2c5a.md"
$overview.Range("B4").Value = "Ready for handoff"
$overview.Range("C4").Value = "Ready for handoff"
$overview.Range("D4").Value = "2016-30-19 02:30:34"

$overview.Range("A5").Value = "f235f8b0-97fc-4c24-b46d-e2250464b4b5.md"
$overview.Range("B5").Value = "Ready for handoff"
$overview.Range("C5").Value = "Ready for handoff"
$overview.Range("D5").Value = "2016-30-19 02:30:34"

$overview.Hyperlinks.Add($overview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/cbb11c9f-9aa9-40a5-b27b-62f0134d2c5a.md", "", "", "cbb11c9f-9aa9-40a5-b27b-62f0134d2c5a.md")
$overview.Hyperlinks.Add($overview.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/f235f8b0-97fc-4c24-b46d-e2250464b4b5.md", "", "", "f235f8b0-97fc-4c24-b46d-e2250464b4b5.md")

# ---------------------------------------------------------------------------
# zh-cn sheet: Source File Name | File Extension | Status | Latest Handoff
# File | Latest Handoff Datetime | Latest Target File | Latest Handback File
# | Latest Handback DateTime | Handoff Reason | Dependency From | Error Detail
# ---------------------------------------------------------------------------

$zhcn.Range("A4").Value = "cbb11c9f-9aa9-40a5-b27b-62f0134d2c5a.md"
$zhcn.Range("B4").Value = ".md"
$zhcn.Range("C4").Value = "Ready for handoff"
$zhcn.Range("D4").Value = "cbb11c9f-9aa9-40a5-b27b-62f0134d2c5a.15fc3705ca78e224c7b59058c2f737e0d764ee99.zh-cn.xlf"
$zhcn.Range("E4").Value = "2016-03-19 02:30:31"
$zhcn.Range("H4").Value = "0001-01-01 00:00:00"
$zhcn.Range("I4").Value = "Include"

$zhcn.Range("A5").Value = "f235f8b0-97fc-4c24-b46d-e2250464b4b5.md"
$zhcn.Range("B5").Value = ".md"
$zhcn.Range("C5").Value = "Ready for handoff"
$zhcn.Range("D5").Value = "f235f8b0-97fc-4c24-b46d-e2250464b4b5.3fed449173b93c0372dee947fdd97e4547de91bf.zh-cn.xlf"
$zhcn.Range("E5").Value = "2016-03-19 02:30:31"
$zhcn.Range("H5").Value = "0001-01-01 00:00:00"
$zhcn.Range("I5").Value = "Include"

$zhcn.Hyperlinks.Add($zhcn.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/cbb11c9f-9aa9-40a5-b27b-62f0134d2c5a.md", "", "", "cbb11c9f-9aa9-40a5-b27b-62f0134d2c5a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("B4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/master/e2e/cbb11c9f-9aa9-40a5-b27b-62f0134d2c5a.md", "", "", ".md")
$zhcn.Hyperlinks.Add($zhcn.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/cbb11c9f-9aa9-40a5-b27b-62f0134d2c5a.15fc3705ca78e224c7b59058c2f737e0d764ee99.zh-cn.xlf", "", "", "cbb11c9f-9aa9-40a5-b27b-62f0134d2c5a.15fc3705ca78e224c7b59058c2f737e0d764ee99.zh-cn.xlf")

$zhcn.Hyperlinks.Add($zhcn.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/f235f8b0-97fc-4c24-b46d-e2250464b4b5.md", "", "", "f235f8b0-97fc-4c24-b46d-e2250464b4b5.md")
$zhcn.Hyperlinks.Add($zhcn.Range("B5"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/master/e2e/f235f8b0-97fc-4c24-b46d-e2250464b4b5.md", "", "", ".md")
$zhcn.Hyperlinks.Add($zhcn.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f235f8b0-97fc-4c24-b46d-e2250464b4b5.3fed449173b93c0372dee947fdd97e4547de91bf.zh-cn.xlf", "", "", "f235f8b0-97fc-4c24-b46d-e2250464b4b5.3fed449173b93c0372dee947fdd97e4547de91bf.zh-cn.xlf")

# ---------------------------------------------------------------------------
# de-de sheet: same shape as zh-cn
# ---------------------------------------------------------------------------

$dede.Range("A4").Value = "cbb11c9f-9aa9-40a5-b27b-62f0134d2c5a.md"
$dede.Range("B4").Value = ".md"
$dede.Range("C4").Value = "Ready for handoff"
$dede.Range("D4").Value = "cbb11c9f-9aa9-40a5-b27b-62f0134d2c5a.15fc3705ca78e224c7b59058c2f737e0d764ee99.de-de.xlf"
$dede.Range("E4").Value = "2016-03-19 02:30:34"
$dede.Range("H4").Value = "0001-01-01 00:00:00"
$dede.Range("I4").Value = "Include"

$dede.Range("A5").Value = "f235f8b0-97fc-4c24-b46d-e2250464b4b5.md"
$dede.Range("B5").Value = ".md"
$dede.Range("C5").Value = "Ready for handoff"
$dede.Range("D5").Value = "f235f8b0-97fc-4c24-b46d-e2250464b4b5.3fed449173b93c0372dee947fdd97e4547de91bf.de-de.xlf"
$dede.Range("E5").Value = "2016-03-19 02:30:34"
$dede.Range("H5").Value = "0001-01-01 00:00:00"
$dede.Range("I5").Value = "Include"

$dede.Hyperlinks.Add($dede.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/cbb11c9f-9aa9-40a5-b27b-62f0134d2c5a.md", "", "", "cbb11c9f-9aa9-40a5-b27b-62f0134d2c5a.md")
$dede.Hyperlinks.Add($dede.Range("B4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/master/e2e/cbb11c9f-9aa9-40a5-b27b-62f0134d2c5a.md", "", "", ".md")
$dede.Hyperlinks.Add($dede.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/cbb11c9f-9aa9-40a5-b27b-62f0134d2c5a.15fc3705ca78e224c7b59058c2f737e0d764ee99.de-de.xlf", "", "", "cbb11c9f-9aa9-40a5-b27b-62f0134d2c5a.15fc3705ca78e224c7b59058c2f737e0d764ee99.de-de.xlf")

$dede.Hyperlinks.Add($dede.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/f235f8b0-97fc-4c24-b46d-e2250464b4b5.md", "", "", "f235f8b0-97fc-4c24-b46d-e2250464b4b5.md")
$dede.Hyperlinks.Add($dede.Range("B5"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/master/e2e/f235f8b0-97fc-4c24-b46d-e2250464b4b5.md", "", "", ".md")
$dede.Hyperlinks.Add($dede.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f235f8b0-97fc-4c24-b46d-e2250464b4b5.3fed449173b93c0372dee947fdd97e4547de91bf.de-de.xlf", "", "", "f235f8b0-97fc-4c24-b46d-e2250464b4b5.3fed449173b93c0372dee947fdd97e4547de91bf.de-de.xlf")
